$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.100.27"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3
$ws.Range("D3").Value = "1.909.66"
$ws.Range("E3").Value = "  -1.27%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'0.7406"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.64%  "

# Row 6
$ws.Range("D6").Value = "'244.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.17%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "'0.3093"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.89%  "

# Row 9
$ws.Range("D9").Value = "'26.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.24%  "

# Row 10
$ws.Range("D10").Value = "'0.06978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.70%  "

# Row 11
$ws.Range("D11").Value = "'0.08080"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "

# Row 12
$ws.Range("D12").Value = "'0.7710"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.35%  "

# Row 13
$ws.Range("D13").Value = "1.909.53"
$ws.Range("E13").Value = "  -1.19%  "

# Row 14
$ws.Range("D14").Value = "'5.331"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.33%  "

# Row 15
$ws.Range("D15").Value = "'92.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.11%  "

# Row 16
$ws.Range("D16").Value = "'14.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.82%  "

# Row 17
$ws.Range("D17").Value = "30.115.54"
$ws.Range("E17").Value = "  -0.73%  "

# Row 18
$ws.Range("D18").Value = "'6.072"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "

# Row 19
$ws.Range("D19").Value = "'0.000007844"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.66%  "

# Row 20
$ws.Range("D20").Value = "'240.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.85%  "

# Row 21
$ws.Range("D21").Value = "2.221.31"
$ws.Range("E21").Value = "  +1.45%  "

# Row 22
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").Value = "'1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "

# Row 24
$ws.Range("D24").Value = "'7.097"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.83%  "

# Row 25
$ws.Range("D25").Value = "'9.403"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.25%  "

# Row 26
$ws.Range("D26").Value = "'167.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.67%  "

# Row 27
$ws.Range("D27").Value = "'18.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "

# Row 28
$ws.Range("D28").Value = "'0.1279"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.58%  "

# Row 29
$ws.Range("D29").Value = "'2.058"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.51%  "

# Row 30
$ws.Range("D30").Value = "'1.551"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.16%  "

# Row 31
$ws.Range("D31").Value = "'1.357"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.77%  "

# Row 32
$ws.Range("D32").Value = "'4.341"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.53%  "

# Row 33
$ws.Range("D33").Value = "'4.084"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.21%  "

# Row 34
$ws.Range("D34").Value = "'1.309"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.18%  "

# Row 35
$ws.Range("D35").Value = "'0.05154"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.38%  "

# Row 36
$ws.Range("D36").Value = "'0.7492"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.70%  "

# Row 37
$ws.Range("D37").Value = "'2.725"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.33%  "

# Row 38
$ws.Range("D38").Value = "'0.01961"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.33%  "

# Row 39
$ws.Range("D39").Value = "'2.796"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.51%  "

# Row 40
$ws.Range("D40").Value = "'6.349"
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.4510"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.46%  "

# Row 42
$ws.Range("D42").Value = "'74.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.31%  "

# Row 43
$ws.Range("E43").Value = "  +0.25%  "

# Row 44
$ws.Range("E44").Value = "  +0.06%  "

# Row 45
$ws.Range("D45").Value = "'0.8407"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.37%  "

# Row 46
$ws.Range("D46").Value = "'7.746"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.70%  "

# Row 49
$ws.Range("D49").Value = "2.077.93"
$ws.Range("E49").Value = "  -0.43%  "

# Row 50
$ws.Range("D50").Value = "'36.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.73%  "

# Row 51
$ws.Range("D51").Value = "'0.1186"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.20%  "

# Row 47 (was EnergySwap -> now Quant)
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'102.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.44%  "

# Row 48 (was Quant -> now EnergySwap)
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.969"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.20%  "
